$wb = $excel.ActiveWorkbook

# --- Sheet1: header row (row 6) height tweak 23.25 -> 20.25 ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("B6:H6").RowHeight = 20.25

# --- Add Sheet3 after Sheet2, becoming the new active sheet ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Add($null, $ws2)

# Column widths (closest reachable values to the 95.25 / 44.875 targets)
$ws3.Columns.Item(1).ColumnWidth = 94.57142857142857
$ws3.Columns.Item(2).ColumnWidth = 44.142857142857146

# Row data. NOTE: rows 5 and 6 are filled out of natural order (6 before 5)
# so that shared-string indices land exactly as in the target workbook
# (idx 76 "simple admittance ... pos_table" must be registered before
#  idx 77 "yang's admittance ... pos_table").
$ws3.Range("A1").Value = "yang's 导纳参数；计算机器人本体斥力；不加IMU；不计算障碍物斥力"
$ws3.Range("B1").Value = "可以跑完；但会碰到第五个障碍物"

$ws3.Range("A2").Value = "简单导纳参数；计算机器人本体斥力；不加IMU；不计算障碍物斥力"
$ws3.Range("B2").Value = "可以跑完；但会碰到第五个障碍物；在终点处有徘徊"

$ws3.Range("A3").Value = "yang's 导纳参数；计算机器人本体斥力；加IMU；不计算障碍物斥力"
$ws3.Range("B3").Value = "四五障碍物间超限辐"

$ws3.Range("A4").Value = "简单导纳参数；计算机器人本体斥力；加IMU；不计算障碍物斥力"
$ws3.Range("B4").Value = "四五障碍物间超限辐"

$ws3.Range("A6").Value = "简单导纳参数；计算机器人本体斥力；加IMU；不计算障碍物斥力，同时将pos_table同步为末端位置"
$ws3.Range("B6").Value = "可以跑完；但会碰到第五个障碍物；在终点处有徘徊"

$ws3.Range("A5").Value = "yang's 导纳参数；计算机器人本体斥力；加IMU；不计算障碍物斥力，同时将pos_table同步为末端位置"
$ws3.Range("B5").Value = "可以跑完；但会碰到第五个障碍物"

# Page setup matches Sheet1/Sheet2 (A4/portrait-style print settings)
$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

# Selection matches the committed workbook: B3 selected on Sheet3
$ws3.Range("B3").Select() | Out-Null
